$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("J1:K1")
$r.Value = "test"
$r.Font.Bold = $true
$r.Interior.ThemeColor = 7
$r.Borders.LineStyle = 1
Write-Output "done"
